$wb = $excel.ActiveWorkbook

# --- Step 1: Reorganize sheets ---
$oldTotal = $wb.Worksheets.Item("总计")
# Duplicate the old '总计' sheet; the duplicate will become the new '总计' (gets a fresh sheetId)
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")
# Rename old sheet (keeps its original sheetId) to the new quarter name
$oldTotal.Name = "2022-Q1"
# Now free to rename the duplicate back to 总计
$newTotal.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- Step 2: Build '2022-Q1' sheet content & layout from the '2021-Q4' template ---
$q1.Cells.Clear()
$q4.Range("A1:H25").Copy($q1.Range("A1"))
$q1.Range("A1").Clear()
$q1.Rows.Item(24).Delete()
$q1.Rows.Item(24).Delete()

# Write fund-holding rows as text (preserve formats such as leading zeros / fixed decimals)
$q1Text = $q1.Range("B2:G23")
$q1Text.NumberFormat = "@"
$q1.Cells.Item(2,2).Value = '012930'
$q1.Cells.Item(2,3).Value = '中庚价值先锋股票'
$q1.Cells.Item(2,4).Value = '54.59'
$q1.Cells.Item(2,5).Value = '94.46'
$q1.Cells.Item(2,6).Value = '6.80'
$q1.Cells.Item(2,7).Value = '3.7121'
$q1.Cells.Item(2,8).Value = 2
$q1.Cells.Item(3,2).Value = '013910'
$q1.Cells.Item(3,3).Value = '兴业兴睿两年持有期混合A'
$q1.Cells.Item(3,4).Value = '79.49'
$q1.Cells.Item(3,5).Value = '39.91'
$q1.Cells.Item(3,6).Value = '1.15'
$q1.Cells.Item(3,7).Value = '0.9141'
$q1.Cells.Item(3,8).Value = 8
$q1.Cells.Item(4,2).Value = '501078'
$q1.Cells.Item(4,3).Value = '广发科创主题 3 年封闭运作灵活配置混合'
$q1.Cells.Item(4,4).Value = '21.07'
$q1.Cells.Item(4,5).Value = '95.22'
$q1.Cells.Item(4,6).Value = '2.88'
$q1.Cells.Item(4,7).Value = '0.6068'
$q1.Cells.Item(4,8).Value = 8
$q1.Cells.Item(5,2).Value = '162720'
$q1.Cells.Item(5,3).Value = '广发创业板两年定期开放混合'
$q1.Cells.Item(5,4).Value = '8.96'
$q1.Cells.Item(5,5).Value = '93.33'
$q1.Cells.Item(5,6).Value = '3.85'
$q1.Cells.Item(5,7).Value = '0.3450'
$q1.Cells.Item(5,8).Value = 7
$q1.Cells.Item(6,2).Value = '013911'
$q1.Cells.Item(6,3).Value = '兴业兴睿两年持有期混合C'
$q1.Cells.Item(6,4).Value = '26.65'
$q1.Cells.Item(6,5).Value = '39.91'
$q1.Cells.Item(6,6).Value = '1.15'
$q1.Cells.Item(6,7).Value = '0.3065'
$q1.Cells.Item(6,8).Value = 8
$q1.Cells.Item(7,2).Value = '009353'
$q1.Cells.Item(7,3).Value = '浙商科技创新一个月滚动持有混合A'
$q1.Cells.Item(7,4).Value = '9.37'
$q1.Cells.Item(7,5).Value = '93.15'
$q1.Cells.Item(7,6).Value = '2.99'
$q1.Cells.Item(7,7).Value = '0.2802'
$q1.Cells.Item(7,8).Value = 10
$q1.Cells.Item(8,2).Value = '166801'
$q1.Cells.Item(8,3).Value = '浙商聚潮新思维混合'
$q1.Cells.Item(8,4).Value = '9.62'
$q1.Cells.Item(8,5).Value = '78.06'
$q1.Cells.Item(8,6).Value = '2.71'
$q1.Cells.Item(8,7).Value = '0.2607'
$q1.Cells.Item(8,8).Value = 10
$q1.Cells.Item(9,2).Value = '168501'
$q1.Cells.Item(9,3).Value = '北信瑞丰产业升级多策略混合'
$q1.Cells.Item(9,4).Value = '4.42'
$q1.Cells.Item(9,5).Value = '94.11'
$q1.Cells.Item(9,6).Value = '3.57'
$q1.Cells.Item(9,7).Value = '0.1578'
$q1.Cells.Item(9,8).Value = 9
$q1.Cells.Item(10,2).Value = '009354'
$q1.Cells.Item(10,3).Value = '浙商科技创新一个月滚动持有混合C'
$q1.Cells.Item(10,4).Value = '3.92'
$q1.Cells.Item(10,5).Value = '93.15'
$q1.Cells.Item(10,6).Value = '2.99'
$q1.Cells.Item(10,7).Value = '0.1172'
$q1.Cells.Item(10,8).Value = 10
$q1.Cells.Item(11,2).Value = '014189'
$q1.Cells.Item(11,3).Value = '南方专精特新混合A'
$q1.Cells.Item(11,4).Value = '4.99'
$q1.Cells.Item(11,5).Value = '30.05'
$q1.Cells.Item(11,6).Value = '1.31'
$q1.Cells.Item(11,7).Value = '0.0654'
$q1.Cells.Item(11,8).Value = 10
$q1.Cells.Item(12,2).Value = '001255'
$q1.Cells.Item(12,3).Value = '长城改革红利灵活配置混合'
$q1.Cells.Item(12,4).Value = '1.29'
$q1.Cells.Item(12,5).Value = '60.02'
$q1.Cells.Item(12,6).Value = '3.27'
$q1.Cells.Item(12,7).Value = '0.0422'
$q1.Cells.Item(12,8).Value = 5
$q1.Cells.Item(13,2).Value = '004266'
$q1.Cells.Item(13,3).Value = '招商沪港深科技创新主题精选灵活配置混合A'
$q1.Cells.Item(13,4).Value = '1.29'
$q1.Cells.Item(13,5).Value = '88.85'
$q1.Cells.Item(13,6).Value = '3.05'
$q1.Cells.Item(13,7).Value = '0.0393'
$q1.Cells.Item(13,8).Value = 7
$q1.Cells.Item(14,2).Value = '014190'
$q1.Cells.Item(14,3).Value = '南方专精特新混合C'
$q1.Cells.Item(14,4).Value = '1.13'
$q1.Cells.Item(14,5).Value = '30.05'
$q1.Cells.Item(14,6).Value = '1.31'
$q1.Cells.Item(14,7).Value = '0.0148'
$q1.Cells.Item(14,8).Value = 10
$q1.Cells.Item(15,2).Value = '011214'
$q1.Cells.Item(15,3).Value = '招商惠润一年定期开放混合型发起式管理人中管理人（MOM）证券投资基金A'
$q1.Cells.Item(15,4).Value = '0.67'
$q1.Cells.Item(15,5).Value = '81.20'
$q1.Cells.Item(15,6).Value = '2.05'
$q1.Cells.Item(15,7).Value = '0.0137'
$q1.Cells.Item(15,8).Value = 9
$q1.Cells.Item(16,2).Value = '009128'
$q1.Cells.Item(16,3).Value = '明亚价值长青混合A'
$q1.Cells.Item(16,4).Value = '0.38'
$q1.Cells.Item(16,5).Value = '49.48'
$q1.Cells.Item(16,6).Value = '2.81'
$q1.Cells.Item(16,7).Value = '0.0107'
$q1.Cells.Item(16,8).Value = 7
$q1.Cells.Item(17,2).Value = '010754'
$q1.Cells.Item(17,3).Value = '招商沪港深科技创新主题精选灵活配置混合C'
$q1.Cells.Item(17,4).Value = '0.28'
$q1.Cells.Item(17,5).Value = '88.85'
$q1.Cells.Item(17,6).Value = '3.05'
$q1.Cells.Item(17,7).Value = '0.0085'
$q1.Cells.Item(17,8).Value = 7
$q1.Cells.Item(18,2).Value = '004521'
$q1.Cells.Item(18,3).Value = '安信工业4.0主题沪港深精选灵活配置混合A'
$q1.Cells.Item(18,4).Value = '0.09'
$q1.Cells.Item(18,5).Value = '85.98'
$q1.Cells.Item(18,6).Value = '5.68'
$q1.Cells.Item(18,7).Value = '0.0051'
$q1.Cells.Item(18,8).Value = 7
$q1.Cells.Item(19,2).Value = '002303'
$q1.Cells.Item(19,3).Value = '金鹰智慧生活灵活配置混合'
$q1.Cells.Item(19,4).Value = '0.11'
$q1.Cells.Item(19,5).Value = '89.88'
$q1.Cells.Item(19,6).Value = '4.33'
$q1.Cells.Item(19,7).Value = '0.0048'
$q1.Cells.Item(19,8).Value = 6
$q1.Cells.Item(20,2).Value = '001866'
$q1.Cells.Item(20,3).Value = '北信瑞丰新成长灵活配置混合'
$q1.Cells.Item(20,4).Value = '0.07'
$q1.Cells.Item(20,5).Value = '94.21'
$q1.Cells.Item(20,6).Value = '5.35'
$q1.Cells.Item(20,7).Value = '0.0037'
$q1.Cells.Item(20,8).Value = 3
$q1.Cells.Item(21,2).Value = '004522'
$q1.Cells.Item(21,3).Value = '安信工业4.0主题沪港深精选灵活配置混合C'
$q1.Cells.Item(21,4).Value = '0.05'
$q1.Cells.Item(21,5).Value = '85.98'
$q1.Cells.Item(21,6).Value = '5.68'
$q1.Cells.Item(21,7).Value = '0.0028'
$q1.Cells.Item(21,8).Value = 7
$q1.Cells.Item(22,2).Value = '011215'
$q1.Cells.Item(22,3).Value = '招商惠润一年定期开放混合型发起式管理人中管理人（MOM）证券投资基金C'
$q1.Cells.Item(22,4).Value = '0.09'
$q1.Cells.Item(22,5).Value = '81.20'
$q1.Cells.Item(22,6).Value = '2.05'
$q1.Cells.Item(22,7).Value = '0.0018'
$q1.Cells.Item(22,8).Value = 9
$q1.Cells.Item(23,2).Value = '009129'
$q1.Cells.Item(23,3).Value = '明亚价值长青混合C'
$q1.Cells.Item(23,4).Value = '0.00'
$q1.Cells.Item(23,5).Value = '49.48'
$q1.Cells.Item(23,6).Value = '2.81'
$q1.Cells.Item(23,8).Value = 7
$q1Text.Style = "Normal"
# G23 is numeric 0 (unlike the other text-formatted G-column cells)
$q1.Cells.Item(23,7).Value = 0

# --- Step 3: Build new '总计' sheet content (insert the 2022-Q1 summary row at the top) ---
$total.Rows.Item(2).Insert()
$total.Range("A3").Copy($total.Range("A2"))
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$totalText = $total.Range("B2:B2")
$totalText.NumberFormat = "@"
$total.Cells.Item(2,2).Value = '2022-Q1'
$totalText.Style = "Normal"
$total.Cells.Item(2,3).Value = 22
$total.Cells.Item(2,4).Value = 6.91

Write-Output "edit complete"
